$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New block of rows appended below the existing table (rows 35-37) ---
# Row 35: label cell (new shared string "pen-level")
$ws.Range("A35").Value = "pen-level"

# Row 36: raw data values
$ws.Range("A36").Value = 1745.4
$ws.Range("B36").Value = 70.1923
$ws.Range("C36").Value = 13801.1
$ws.Range("D36").Value = 58.3972
$ws.Range("E36").Value = 89.0394
$ws.Range("F36").Value = 0.231328
$ws.Range("G36").Value = 34.49
$ws.Range("H36").Value = 0.834787
$ws.Range("I36").Value = 10.7068
$ws.Range("J36").Value = 0.0794539
$ws.Range("K36").Value = 89.32
$ws.Range("L36").Value = 5.00682
$ws.Range("M36").Value = 12707.7
$ws.Range("N36").Value = 185.508
$ws.Range("O36").Value = 0.696121
$ws.Range("P36").Value = 0.0395163

# Row 37: scaled values (0.196 * row 36), same pattern used by every other
# block in the sheet - first cell is its own formula, the rest share it.
# (Filled in before the row-36 number formats below, so the fill doesn't
# drag row 36's number format down into row 37.)
$ws.Range("A37").Formula = "=0.196 * A36"
$ws.Range("B37:P37").Formula = "=0.196 * B36"

# A handful of row-36 cells carry an explicit "#,##0" number format
# (matches the new cellXfs entry with numFmtId 3 in the edited file).
$ws.Range("B36").NumberFormat = "#,##0"
$ws.Range("D36").NumberFormat = "#,##0"
$ws.Range("E36").NumberFormat = "#,##0"
$ws.Range("I36").NumberFormat = "#,##0"
$ws.Range("L36").NumberFormat = "#,##0"
$ws.Range("N36").NumberFormat = "#,##0"

# --- View state: last selected cell in the refreshed sheet ---
$ws.Range("N37").Select()

# --- Page setup: paper size / orientation were stamped onto the sheet ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
